# Auto-applied data refresh for cryptos.xlsx (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.583.68"
$ws.Range("E2").Value = "  -2.53%  "
$ws.Range("D3").Value = "2.344.04"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'500.96"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").Value = "'130.04"
$ws.Range("E6").Value = "  -3.66%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'0.535"
$ws.Range("E8").Value = "  -3.65%  "
$ws.Range("D9").Value = "2.357.66"
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").Value = "'0.0949"
$ws.Range("E10").Value = "  -3.77%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "'4.74"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "'0.318"
$ws.Range("E13").Value = "  -5.37%  "
$ws.Range("D14").Value = "2.764.60"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "55.567.26"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'21.54"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("D18").Value = "2.370.80"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'9.77"
$ws.Range("E19").Value = "  -4.35%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'307.24"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "'3.99"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "'6.18"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'64.83"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").Value = "'0.145"
$ws.Range("E27").Value = "  -4.66%  "
$ws.Range("D28").Value = "'7.14"
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("D29").Value = "'169.12"
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("D30").Value = "0.0₃0702"
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("D31").Value = "'1.63"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'5.73"
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("D35").Value = "'1.07"
$ws.Range("E35").Value = "  -5.95%  "
$ws.Range("D36").Value = "'17.59"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("D37").Value = "'1.17"
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("D38").Value = "'0.849"
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("D39").Value = "'3.65"
$ws.Range("E39").Value = "  -5.91%  "
$ws.Range("D40").Value = "'36.17"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.373"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.38"
$ws.Range("E42").Value = "  -4.49%  "
$ws.Range("D43").Value = "'3.36"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").Value = "'4.81"
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("D45").Value = "'124.04"
$ws.Range("E45").Value = "  -6.16%  "
$ws.Range("D46").Value = "'0.552"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("D47").Value = "'0.0888"
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("D48").Value = "'240.34"
$ws.Range("E48").Value = "  -4.05%  "
$ws.Range("D49").Value = "'0.0477"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").Value = "'16.80"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("D51").Value = "'0.0205"
$ws.Range("E51").Value = "  -2.72%  "
